# Applies the "Github Link" hyperlink-fields edit:
#  - bolds the paragraph mark of the "Github Link:  " paragraph
#  - appends a HYPERLINK field (complex field code) for the repo root URL
#    to that same paragraph
#  - inserts a brand new paragraph (indented) containing a second
#    HYPERLINK field for the .../tree/main/A2 URL
#  - moves the (hidden) _GoBack bookmark from the "Github Link:" paragraph
#    down to the last (now-5th) paragraph of that block
#
# iron_native's InsertXML always starts a fresh <w:p>, and it silently
# drops <w:rStyle> from inline rPr, so:
#   * each field-run-sequence is inserted as its own paragraph, then the
#     paragraph break in front of it is deleted when it needs to merge
#     into the previous paragraph (hyperlink #1 only; hyperlink #2 stays
#     in its own new paragraph).
#   * the "https://..." display-text run's Hyperlink character style is
#     applied afterwards via Range.Style = "Hyperlink" (which correctly
#     serialises to <w:rStyle w:val="7"/>).
# Bookmarks.Add() on a Range that contains no real characters (a
# collapsed point, or a range spanning only the paragraph mark) comes
# back mis-anchored, so the bookmark is re-created around a throwaway
# character that is deleted again right after.

$d = $word.ActiveDocument

function New-FieldParagraphXml([string]$url, [bool]$withIndent) {
    $runRPr = '<w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="22"/><w:lang w:val="en-IN"/></w:rPr>'
    $pPrRPr = $runRPr
    if ($withIndent) {
        $pPr = '<w:pPr><w:spacing w:after="0"/><w:ind w:left="720" w:leftChars="0" w:firstLine="720" w:firstLineChars="0"/>' + $pPrRPr + '</w:pPr>'
    } else {
        $pPr = '<w:pPr><w:spacing w:after="0"/>' + $pPrRPr + '</w:pPr>'
    }
    $runs = ''
    $runs += '<w:r>' + $runRPr + '<w:fldChar w:fldCharType="begin"/></w:r>'
    $runs += '<w:r>' + $runRPr + '<w:instrText xml:space="preserve"> HYPERLINK "' + $url + '" </w:instrText></w:r>'
    $runs += '<w:r>' + $runRPr + '<w:fldChar w:fldCharType="separate"/></w:r>'
    $runs += '<w:r>' + $runRPr + '<w:t>' + $url + '</w:t></w:r>'
    $runs += '<w:r>' + $runRPr + '<w:fldChar w:fldCharType="end"/></w:r>'
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pPr + $runs + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- locate the "Github Link:  " paragraph -----------------------------
$ghPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Github Link:*") {
        $ghPara = $d.Paragraphs($i)
        break
    }
}

$ghEnd = $ghPara.Range.End            # position right after the paragraph mark
$insPoint = $ghEnd - 1                # position right before the paragraph mark

# --- 1) bold the paragraph mark (+ keep font) of the Github Link para --
$ghPara.Range.Font.Bold = 1
$ghPara.Range.Font.Name = "Cambria"

# --- 2) insert hyperlink #1 field as a brand-new paragraph, then merge
#        it into the Github Link paragraph by deleting the separating
#        paragraph mark -----------------------------------------------
$xml1 = New-FieldParagraphXml "https://github.com/Ganesh-Chavhan/HPC_LAB" $false
$pt1 = $d.Range($insPoint, $insPoint)
$pt1.InsertXML($xml1)

$mergeMark = $d.Range($insPoint, $insPoint + 1)
$mergeMark.Delete()

# apply the Hyperlink character style to the display-text run of field #1
$find1 = $d.Range(0, 0)
$find1.Find.Execute("https://github.com/Ganesh-Chavhan/HPC_LAB", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$find1.Style = "Hyperlink"

# --- 3) insert hyperlink #2 field as its own new (indented) paragraph,
#        right after the Github Link paragraph -------------------------
$ghParaEnd = $ghPara.Range.End
$xml2 = New-FieldParagraphXml "https://github.com/Ganesh-Chavhan/HPC_LAB/tree/main/A2" $true
$pt2 = $d.Range($ghParaEnd - 1, $ghParaEnd - 1)
$pt2.InsertXML($xml2)

$find2 = $d.Range(0, 0)
$find2.Find.Execute("https://github.com/Ganesh-Chavhan/HPC_LAB/tree/main/A2", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$find2.Style = "Hyperlink"

# --- 4) move the _GoBack bookmark from the Github Link paragraph to the
#        last paragraph of the block (the trailing sz-28 empty one) -----
$oldBm = $null
foreach ($nm in @("_GoBack")) {
    $b = $d.Bookmarks($nm)
    if ($b -ne $null) { $oldBm = $b }
}
if ($oldBm -ne $null) { $oldBm.Delete() }

# find the trailing empty "sz 28" paragraph (last paragraph of the doc body)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lpStart = $lastPara.Range.Start
$dummy = $d.Range($lpStart, $lpStart)
$dummy.InsertAfter("X")
$lastPara2 = $d.Paragraphs($d.Paragraphs.Count)
$bmRange = $d.Range($lastPara2.Range.Start, $lastPara2.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$dummyChar = $d.Range($lastPara2.Range.Start, $lastPara2.Range.Start + 1)
$dummyChar.Delete()

Write-Output "done"
